# Applies the commit: removes stray <w:proofErr> spell-check markers,
# merges the "Parsing" run with its follow-up run, and replaces the
# trailing "16h38 - " line with a new "Samedi 5" heading, a "16h45 - "
# line (built from several runs), and a blank paragraph.

$d = $word.ActiveDocument

$wNs  = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$w14Ns = "http://schemas.microsoft.com/office/word/2010/wordml"
$nsAttrs = 'xmlns:w="' + $wNs + '" xmlns:w14="' + $w14Ns + '"'

function Get-ParagraphByText($doc, [string]$needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# --- 1. "Organisation des <<dataset>>" paragraph: drop proofErr around "dataset" and around the red "10m" ---
$xmlOrganisation = @'
<w:p {NS} w14:paraId="6906CE26" w14:textId="74682720" w:rsidR="000F44D5" w:rsidRDefault="00744BCB" w:rsidP="00744BCB"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Organisation des « </w:t></w:r><w:r w:rsidR="00B73049"><w:t>dataset</w:t></w:r><w:r w:rsidR="00D523EE"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">» – </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>10m</w:t></w:r><w:r w:rsidR="00C92A21"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C92A21"><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>10m</w:t></w:r></w:p>
'@
$xmlOrganisation = $xmlOrganisation.Replace("{NS}", $nsAttrs).Trim()
$p = Get-ParagraphByText $d "dataset"
$p.Range.InsertXML($xmlOrganisation)

# --- 2. "Table de correspondance dans Indices.txt" paragraph: drop proofErr around the red "5m" ---
$xmlIndices = @'
<w:p {NS} w14:paraId="1C31E794" w14:textId="7A715D88" w:rsidR="00744BCB" w:rsidRPr="0037293E" w:rsidRDefault="00744BCB" w:rsidP="00744BCB"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Table de correspondance dans Indices.txt – </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>5m</w:t></w:r><w:r w:rsidR="000515B5"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="000515B5"><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>5m</w:t></w:r></w:p>
'@
$xmlIndices = $xmlIndices.Replace("{NS}", $nsAttrs).Trim()
$p = Get-ParagraphByText $d "Table de correspondance dans Indices.txt"
$p.Range.InsertXML($xmlIndices)

# --- 3. "Parsing de la table de correspondance" paragraph: drop proofErr around "Parsing" and
#        merge its run with the following " de la table de correspondance – " run ---
$xmlParsing = @'
<w:p {NS} w14:paraId="0897DE09" w14:textId="78313BB5" w:rsidR="00744BCB" w:rsidRDefault="00744BCB" w:rsidP="00744BCB"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Parsing de la table de correspondance – </w:t></w:r><w:r w:rsidRPr="00744BCB"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>10m</w:t></w:r><w:r w:rsidR="00C92A21"><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> 1m</w:t></w:r></w:p>
'@
$xmlParsing = $xmlParsing.Replace("{NS}", $nsAttrs).Trim()
$p = Get-ParagraphByText $d "Parsing"
$p.Range.InsertXML($xmlParsing)

# --- 4. "Gestion du chronomètre" paragraph: drop proofErr around the red "3h" ---
$xmlChrono = @'
<w:p {NS} w14:paraId="432755BB" w14:textId="1577810B" w:rsidR="00744BCB" w:rsidRDefault="00744BCB" w:rsidP="00744BCB"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Gestion du chronomètre – </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>3h</w:t></w:r><w:r w:rsidR="00C92A21"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C92A21"><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>3h</w:t></w:r></w:p>
'@
$xmlChrono = $xmlChrono.Replace("{NS}", $nsAttrs).Trim()
$p = Get-ParagraphByText $d "Gestion du chronomètre"
$p.Range.InsertXML($xmlChrono)

# --- 5. "Nettoyage de code" paragraph: drop proofErr around the red "3h" ---
$xmlNettoyage = @'
<w:p {NS} w14:paraId="60808E99" w14:textId="42852956" w:rsidR="00335BF5" w:rsidRPr="005739D3" w:rsidRDefault="00335BF5" w:rsidP="00744BCB"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r w:rsidRPr="00335BF5"><w:t>Nettoyage</w:t></w:r><w:r><w:t xml:space="preserve"> de code – </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">3h </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>3h</w:t></w:r></w:p>
'@
$xmlNettoyage = $xmlNettoyage.Replace("{NS}", $nsAttrs).Trim()
$p = Get-ParagraphByText $d "Nettoyage"
$p.Range.InsertXML($xmlNettoyage)

# --- 6. Replace the trailing "16h38 - " paragraph with "Samedi 5" / "16h45 - " / blank paragraph ---
$xmlTail = @'
<w:p {NS}><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Samedi 5</w:t></w:r></w:p><w:p {NS}><w:r><w:t>1</w:t></w:r><w:r><w:t>6</w:t></w:r><w:r><w:t>h</w:t></w:r><w:r><w:t>45</w:t></w:r><w:r><w:t xml:space="preserve"> - </w:t></w:r></w:p><w:p {NS}/>
'@
$xmlTail = $xmlTail.Replace("{NS}", $nsAttrs).Trim()
$p = Get-ParagraphByText $d "16h38"
$p.Range.InsertXML($xmlTail)
